# Scheduled runner update: refresh Kraken market-board price/profit figures
# across the Leve tracking sheets (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 9
$ws.Range("H9").Value = 285.5
$ws.Range("I9").Value = 180.66667
$ws.Range("K9").Value = 180.66667
$ws.Range("M9").Value = -11.66667000000001

# row 15
$ws.Range("H15").Value = 330.14816
$ws.Range("I15").Value = 330.14816
$ws.Range("K15").Value = 990.4444800000001
$ws.Range("M15").Value = -821.4444800000001

# row 32
$ws.Range("H32").Value = 9164.223
$ws.Range("I32").Value = 6239
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 6239
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -5913
$ws.Range("N32").Value = -10652

# row 40
$ws.Range("H40").Value = 7136.625
$ws.Range("J40").Value = 8749.833000000001
$ws.Range("L40").Value = 8749.833000000001
$ws.Range("N40").Value = -9099.833000000001

# row 51
$ws.Range("H51").Value = 7250
$ws.Range("J51").Value = 7250
$ws.Range("L51").Value = 7250
$ws.Range("N51").Value = -8218

# row 80
$ws.Range("H80").Value = 950
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -4996

# row 83
$ws.Range("H83").Value = 950
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -3108
$ws.Range("N83").Value = -18984

# row 92
$ws.Range("H92").Value = 641.1667
$ws.Range("I92").Value = 461.75
$ws.Range("K92").Value = 461.75
$ws.Range("M92").Value = 786.25

# row 97
$ws.Range("H97").Value = 7344
$ws.Range("J97").Value = 7344
$ws.Range("L97").Value = 22032
$ws.Range("N97").Value = -23024

# row 107
$ws.Range("H107").Value = 1077.4286
$ws.Range("I107").Value = 548.75
$ws.Range("J107").Value = 4249.5
$ws.Range("K107").Value = 548.75
$ws.Range("L107").Value = 4249.5
$ws.Range("M107").Value = 1371.25
$ws.Range("N107").Value = -8089.5

# row 111
$ws.Range("H111").Value = 4803.3335
$ws.Range("I111").Value = 3391.2222
$ws.Range("J111").Value = 6921.5
$ws.Range("K111").Value = 10173.6666
$ws.Range("L111").Value = 20764.5
$ws.Range("M111").Value = -7106.6666
$ws.Range("N111").Value = -26898.5

# row 132
$ws.Range("H132").Value = 9499.571
$ws.Range("I132").Value = 7749.25
$ws.Range("K132").Value = 23247.75
$ws.Range("M132").Value = -20717.75

# row 137
$ws.Range("H137").Value = 3867.077
$ws.Range("I137").Value = 3666.111
$ws.Range("K137").Value = 10998.333
$ws.Range("M137").Value = -8448.332999999999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 13249.5625
$ws.Range("I32").Value = 11001.083
$ws.Range("K32").Value = 11001.083
$ws.Range("M32").Value = -10714.083

# row 45
$ws.Range("H45").Value = 2879.8
$ws.Range("I45").Value = 2754.1667
$ws.Range("K45").Value = 2754.1667
$ws.Range("M45").Value = -2377.1667

# row 61
$ws.Range("H61").Value = 2379.6
$ws.Range("I61").Value = 2379.6
$ws.Range("K61").Value = 2379.6
$ws.Range("M61").Value = -2167.6

# row 97
$ws.Range("H97").Value = 384.75
$ws.Range("I97").Value = 384.75
$ws.Range("K97").Value = 384.75
$ws.Range("M97").Value = 111.25

# row 132
$ws.Range("H132").Value = 5675.6924
$ws.Range("I132").Value = 5675.6924
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17027.0772
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14497.0772
$ws.Range("N132").ClearContents() | Out-Null

# row 136
$ws.Range("H136").Value = 2379.6
$ws.Range("I136").Value = 2379.6
$ws.Range("K136").Value = 7138.799999999999
$ws.Range("M136").Value = -4588.799999999999

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 238.11111
$ws.Range("I22").Value = 125
$ws.Range("K22").Value = 125
$ws.Range("M22").Value = 48

# row 86
$ws.Range("H86").Value = 4770.727
$ws.Range("I86").Value = 2095.8
$ws.Range("J86").Value = 6999.8335
$ws.Range("K86").Value = 2095.8
$ws.Range("L86").Value = 6999.8335
$ws.Range("M86").Value = -972.8000000000002
$ws.Range("N86").Value = -9245.833500000001

# row 89
$ws.Range("H89").Value = 4770.727
$ws.Range("I89").Value = 2095.8
$ws.Range("J89").Value = 6999.8335
$ws.Range("K89").Value = 10479
$ws.Range("L89").Value = 34999.1675
$ws.Range("M89").Value = -4863
$ws.Range("N89").Value = -46231.1675

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 2250
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -4248

# row 65
$ws.Range("H65").Value = 2250
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -21240

# row 80
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246

# row 83
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232

# row 107
$ws.Range("H107").Value = 397.6842
$ws.Range("I107").Value = 314.64285
$ws.Range("K107").Value = 314.64285
$ws.Range("M107").Value = 1605.35715

# row 132
$ws.Range("H132").Value = 3000.818
$ws.Range("I132").Value = 2894.3333
$ws.Range("K132").Value = 8682.999899999999
$ws.Range("M132").Value = -6152.999899999999

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 33
$ws.Range("H33").Value = 67.55556
$ws.Range("I33").Value = 63.375
$ws.Range("J33").Value = 101
$ws.Range("K33").Value = 380.25
$ws.Range("L33").Value = 606
$ws.Range("M33").Value = -97.25
$ws.Range("N33").Value = -1172

# row 38
$ws.Range("H38").Value = 127.625
$ws.Range("I38").Value = 33.285713
$ws.Range("K38").Value = 99.857139
$ws.Range("M38").Value = 247.142861

# row 40
$ws.Range("H40").Value = 256.7647

# row 70
$ws.Range("H70").Value = 7597
$ws.Range("I70").Value = 7597
$ws.Range("K70").Value = 22791
$ws.Range("M70").Value = -22476

# row 73
$ws.Range("H73").Value = 7597
$ws.Range("I73").Value = 7597
$ws.Range("K73").Value = 22791
$ws.Range("M73").Value = -21699

# row 117
$ws.Range("H117").Value = 1557.0769
$ws.Range("I117").Value = 639.6
$ws.Range("J117").Value = 2130.5
$ws.Range("K117").Value = 1918.8
$ws.Range("L117").Value = 6391.5
$ws.Range("M117").Value = 1523.2
$ws.Range("N117").Value = -13275.5

# row 129
$ws.Range("H129").Value = 1366
$ws.Range("I129").Value = 1400
$ws.Range("K129").Value = 4200
$ws.Range("M129").Value = 800

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents() | Out-Null

# row 58
$ws.Range("H58").Value = 32501
$ws.Range("J58").Value = 32501
$ws.Range("L58").Value = 32501
$ws.Range("N58").Value = -33055

# row 102
$ws.Range("H102").Value = 3828
$ws.Range("I102").Value = 3828
$ws.Range("K102").Value = 3828
$ws.Range("M102").Value = -2206

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 20
$ws.Range("H20").Value = 70000
$ws.Range("J20").Value = 70000
$ws.Range("L20").Value = 70000
$ws.Range("N20").Value = -70452

# row 57
$ws.Range("H57").Value = 31822
$ws.Range("I57").Value = 1299
$ws.Range("J57").Value = 47083.5
$ws.Range("K57").Value = 1299
$ws.Range("L57").Value = 47083.5
$ws.Range("M57").Value = -733
$ws.Range("N57").Value = -48215.5

# row 59
$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -41308

# row 68
$ws.Range("H68").Value = 2890.6365
$ws.Range("I68").Value = 2829.7
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 2829.7
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2080.7
$ws.Range("N68").Value = -4998

# row 71
$ws.Range("H71").Value = 2890.6365
$ws.Range("I71").Value = 2829.7
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 14148.5
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -10404.5
$ws.Range("N71").Value = -24988

# row 110
$ws.Range("H110").Value = 44997.5
$ws.Range("J110").Value = 44997.5
$ws.Range("L110").Value = 44997.5
$ws.Range("N110").Value = -53177.5

# row 122
$ws.Range("H122").Value = 6733
$ws.Range("I122").Value = 6946.5
$ws.Range("K122").Value = 20839.5
$ws.Range("M122").Value = -18389.5
